$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 374.1613
$ws.Range("I92").Value = 325.04
$ws.Range("J92").Value = 578.8333
$ws.Range("K92").Value = 325.04
$ws.Range("L92").Value = 578.8333
$ws.Range("M92").Value = 922.96
$ws.Range("N92").Value = -3074.8333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("H6").Value = 18967.334
$ws.Range("I6").Value = 55002
$ws.Range("J6").Value = 950
$ws.Range("K6").Value = 55002
$ws.Range("L6").Value = 950
$ws.Range("M6").Value = -54829
$ws.Range("N6").Value = -1296
$ws.Range("H9").Value = 29333.334
$ws.Range("J9").Value = 29333.334
$ws.Range("L9").Value = 29333.334
$ws.Range("N9").Value = -29673.334
$ws.Range("H20").Value = 29333.334
$ws.Range("J20").Value = 29333.334
$ws.Range("L20").Value = 29333.334
$ws.Range("N20").Value = -29873.334
$ws.Range("H23").Value = 21876.188
$ws.Range("J23").Value = 10769.77
$ws.Range("L23").Value = 10769.77
$ws.Range("N23").Value = -11287.77
$ws.Range("H32").Value = 13335597
$ws.Range("I32").Value = 15153645
$ws.Range("J32").Value = 3249
$ws.Range("K32").Value = 15153645
$ws.Range("L32").Value = 3249
$ws.Range("M32").Value = -15153358
$ws.Range("N32").Value = -3823
$ws.Range("H37").Value = 17017
$ws.Range("I37").Value = 9034
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 9034
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -8761
$ws.Range("N37").Value = -25546
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -25976
$ws.Range("H55").Value = 101
$ws.Range("I55").Value = 101
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 101
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 214
$ws.Range("N55").Value = ""
$ws.Range("H61").Value = 20834848
$ws.Range("I61").Value = 25642460
$ws.Range("J61").Value = 1863.3334
$ws.Range("K61").Value = 25642460
$ws.Range("L61").Value = 1863.3334
$ws.Range("M61").Value = -25642248
$ws.Range("N61").Value = -2287.3334
$ws.Range("H74").Value = 2274.7144
$ws.Range("I74").Value = 583.0540999999999
$ws.Range("J74").Value = 4682.077
$ws.Range("K74").Value = 583.0540999999999
$ws.Range("L74").Value = 4682.077
$ws.Range("M74").Value = 290.9459000000001
$ws.Range("N74").Value = -6430.077
$ws.Range("H77").Value = 2274.7144
$ws.Range("I77").Value = 583.0540999999999
$ws.Range("J77").Value = 4682.077
$ws.Range("K77").Value = 2915.2705
$ws.Range("L77").Value = 23410.385
$ws.Range("M77").Value = 1452.7295
$ws.Range("N77").Value = -32146.385
$ws.Range("H80").Value = 36833.332
$ws.Range("I80").Value = 10500
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -9502
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 36833.332
$ws.Range("I83").Value = 10500
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -26508
$ws.Range("N83").Value = -159984
$ws.Range("H122").Value = 1282.1666
$ws.Range("I122").Value = 839.775
$ws.Range("J122").Value = 2546.1428
$ws.Range("K122").Value = 2519.325
$ws.Range("L122").Value = 7638.428400000001
$ws.Range("M122").Value = -69.32499999999982
$ws.Range("N122").Value = -12538.4284
$ws.Range("H136").Value = 20834848
$ws.Range("I136").Value = 25642460
$ws.Range("J136").Value = 1863.3334
$ws.Range("K136").Value = 76927380
$ws.Range("L136").Value = 5590.0002
$ws.Range("M136").Value = -76924830
$ws.Range("N136").Value = -10690.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1307.2727
$ws.Range("I94").Value = 1211.1111
$ws.Range("J94").Value = 1740
$ws.Range("K94").Value = 1211.1111
$ws.Range("L94").Value = 1740
$ws.Range("M94").Value = -760.1111000000001
$ws.Range("N94").Value = -2642
$ws.Range("H134").Value = 1618.1714
$ws.Range("I134").Value = 1372.2069
$ws.Range("J134").Value = 2807
$ws.Range("K134").Value = 4116.620699999999
$ws.Range("L134").Value = 8421
$ws.Range("M134").Value = -1581.620699999999
$ws.Range("N134").Value = -13491
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 643.12
$ws.Range("I5").Value = 379.5
$ws.Range("J5").Value = 978.63635
$ws.Range("K5").Value = 1138.5
$ws.Range("L5").Value = 2935.90905
$ws.Range("M5").Value = -1026.5
$ws.Range("N5").Value = -3159.90905
$ws.Range("H12").Value = 37.826088
$ws.Range("J12").Value = 42.8
$ws.Range("L12").Value = 128.4
$ws.Range("N12").Value = -474.4
$ws.Range("H97").Value = 362.5
$ws.Range("I97").Value = 175
$ws.Range("J97").Value = 550
$ws.Range("K97").Value = 525
$ws.Range("L97").Value = 1650
$ws.Range("M97").Value = -29
$ws.Range("N97").Value = -2642
$ws.Range("H98").Value = 520.7143
$ws.Range("I98").Value = 392.30768
$ws.Range("J98").Value = 729.375
$ws.Range("K98").Value = 1176.92304
$ws.Range("L98").Value = 2188.125
$ws.Range("M98").Value = 321.0769599999999
$ws.Range("N98").Value = -5184.125
$ws.Range("H118").Value = 646.75
$ws.Range("I118").Value = 646.75
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1940.25
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -697.25
$ws.Range("N118").Value = ""
$ws.Range("H135").Value = 643.12
$ws.Range("I135").Value = 379.5
$ws.Range("J135").Value = 978.63635
$ws.Range("K135").Value = 3415.5
$ws.Range("L135").Value = 8807.727150000001
$ws.Range("M135").Value = -880.5
$ws.Range("N135").Value = -13877.72715
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26372
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81864
$ws.Range("H80").Value = 3001.25
$ws.Range("I80").Value = 2302.5
$ws.Range("J80").Value = 3700
$ws.Range("K80").Value = 2302.5
$ws.Range("L80").Value = 3700
$ws.Range("M80").Value = -1304.5
$ws.Range("N80").Value = -5696
$ws.Range("H83").Value = 3001.25
$ws.Range("I83").Value = 2302.5
$ws.Range("J83").Value = 3700
$ws.Range("K83").Value = 11512.5
$ws.Range("L83").Value = 18500
$ws.Range("M83").Value = -6520.5
$ws.Range("N83").Value = -28484
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2327.5454
$ws.Range("I46").Value = 3220.2
$ws.Range("J46").Value = 1583.6666
$ws.Range("K46").Value = 3220.2
$ws.Range("L46").Value = 1583.6666
$ws.Range("M46").Value = -3032.2
$ws.Range("N46").Value = -1959.6666
$ws.Range("H55").Value = 164.40909
$ws.Range("I55").Value = 129
$ws.Range("J55").Value = 199.81818
$ws.Range("K55").Value = 129
$ws.Range("L55").Value = 199.81818
$ws.Range("M55").Value = 44
$ws.Range("N55").Value = -545.81818
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H82").Value = 1233.3334
$ws.Range("I82").Value = 1050
$ws.Range("J82").Value = 1600
$ws.Range("K82").Value = 1050
$ws.Range("L82").Value = 1600
$ws.Range("M82").Value = -689
$ws.Range("N82").Value = -2322
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H85").Value = 1233.3334
$ws.Range("I85").Value = 1050
$ws.Range("J85").Value = 1600
$ws.Range("K85").Value = 1050
$ws.Range("L85").Value = 1600
$ws.Range("M85").Value = 198
$ws.Range("N85").Value = -4096
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = ""
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = ""
$ws.Range("H122").Value = 104528
$ws.Range("I122").Value = 186151.33
$ws.Range("J122").Value = 6580
$ws.Range("K122").Value = 558453.99
$ws.Range("L122").Value = 19740
$ws.Range("M122").Value = -556003.99
$ws.Range("N122").Value = -24640
$ws.Range("H136").Value = 4694.5137
$ws.Range("I136").Value = 1795.0869
$ws.Range("J136").Value = 9457.857
$ws.Range("K136").Value = 5385.2607
$ws.Range("L136").Value = 28373.571
$ws.Range("M136").Value = -2835.2607
$ws.Range("N136").Value = -33473.571
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1750
$ws.Range("I81").Value = 1666.6666
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3333.3332
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -2272.3332
$ws.Range("N81").Value = -6122
$ws.Range("H82").Value = 30150.25
$ws.Range("J82").Value = 30150.25
$ws.Range("L82").Value = 30150.25
$ws.Range("N82").Value = -30916.25
$ws.Range("H84").Value = 1750
$ws.Range("I84").Value = 1666.6666
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 16666.666
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -11362.666
$ws.Range("N84").Value = -30608
$ws.Range("H85").Value = 30150.25
$ws.Range("J85").Value = 30150.25
$ws.Range("L85").Value = 30150.25
$ws.Range("N85").Value = -32802.25
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H88").Value = 37189
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 37189
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 37189
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -38001
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H91").Value = 37189
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 37189
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 37189
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -39997
$ws.Range("H122").Value = 14359250
$ws.Range("I122").Value = 22335078
$ws.Range("J122").Value = 2760
$ws.Range("K122").Value = 67005234
$ws.Range("L122").Value = 8280
$ws.Range("M122").Value = -67002784
$ws.Range("N122").Value = -13180
